$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 3980
$ws.Range("I20").Value = 3980
$ws.Range("K20").Value = 3980
$ws.Range("M20").Value = -3750
$ws.Range("H35").Value = 3980
$ws.Range("I35").Value = 3980
$ws.Range("K35").Value = 3980
$ws.Range("M35").Value = -3601
$ws.Range("H44").Value = 42968.332
$ws.Range("J44").Value = 42968.332
$ws.Range("L44").Value = 42968.332
$ws.Range("N44").Value = -43892.332
$ws.Range("H111").Value = 1292.9231
$ws.Range("I111").Value = 1072
$ws.Range("K111").Value = 3216
$ws.Range("M111").Value = -149
$ws.Range("H116").Value = 391388.47
$ws.Range("J116").Value = 10253.643
$ws.Range("L116").Value = 10253.643
$ws.Range("N116").Value = -17137.643
$ws.Range("H129").Value = 942.4666999999999
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 958.3953
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 2875.1859
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -12875.1859
$ws.Range("H132").Value = 691740.0600000001
$ws.Range("I132").Value = 401400.25
$ws.Range("J132").Value = 2506364
$ws.Range("K132").Value = 1204200.75
$ws.Range("L132").Value = 7519092
$ws.Range("M132").Value = -1201670.75
$ws.Range("N132").Value = -7524152
$ws.Range("H137").Value = 3361.0667
$ws.Range("I137").Value = 1802.25
$ws.Range("J137").Value = 5142.5713
$ws.Range("K137").Value = 5406.75
$ws.Range("L137").Value = 15427.7139
$ws.Range("M137").Value = -2856.75
$ws.Range("N137").Value = -20527.7139
$ws.Range("H138").Value = 1771.23
$ws.Range("I138").Value = 665.3
$ws.Range("J138").Value = 2508.5166
$ws.Range("K138").Value = 1995.9
$ws.Range("L138").Value = 7525.5498
$ws.Range("M138").Value = 3144.1
$ws.Range("N138").Value = -17805.5498
$ws.Range("H141").Value = 5336.354
$ws.Range("I141").Value = 5571.25
$ws.Range("J141").Value = 2752.5
$ws.Range("K141").Value = 16713.75
$ws.Range("L141").Value = 8257.5
$ws.Range("M141").Value = -11533.75
$ws.Range("N141").Value = -18617.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4575.5205
$ws.Range("I32").Value = 3508.5085
$ws.Range("J32").Value = 9072.214
$ws.Range("K32").Value = 3508.5085
$ws.Range("L32").Value = 9072.214
$ws.Range("M32").Value = -3221.5085
$ws.Range("N32").Value = -9646.214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 822.7646999999999
$ws.Range("I94").Value = 876.9286
$ws.Range("J94").Value = 570
$ws.Range("K94").Value = 876.9286
$ws.Range("L94").Value = 570
$ws.Range("M94").Value = -425.9286
$ws.Range("N94").Value = -1472
$ws.Range("H134").Value = 1840.7838
$ws.Range("I134").Value = 1337.3704
$ws.Range("J134").Value = 3200
$ws.Range("K134").Value = 4012.1112
$ws.Range("L134").Value = 9600
$ws.Range("M134").Value = -1477.1112
$ws.Range("N134").Value = -14670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13890304
$ws.Range("H31").Value = 14708573
$ws.Range("I31").Value = 1560.5238
$ws.Range("J31").Value = 38466056
$ws.Range("K31").Value = 1560.5238
$ws.Range("L31").Value = 38466056
$ws.Range("M31").Value = -1265.5238
$ws.Range("N31").Value = -38466646
$ws.Range("H34").Value = 14708573
$ws.Range("I34").Value = 1560.5238
$ws.Range("J34").Value = 38466056
$ws.Range("K34").Value = 1560.5238
$ws.Range("L34").Value = 38466056
$ws.Range("M34").Value = -1358.5238
$ws.Range("N34").Value = -38466460
$ws.Range("H58").Value = 1503.086
$ws.Range("I58").Value = 1268.6666
$ws.Range("J58").Value = 3691
$ws.Range("K58").Value = 1268.6666
$ws.Range("L58").Value = 3691
$ws.Range("M58").Value = -1065.6666
$ws.Range("N58").Value = -4097
$ws.Range("H113").Value = 13890304
$ws.Range("H132").Value = 1654.775
$ws.Range("I132").Value = 888.0294
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 2664.0882
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -134.0882000000001
$ws.Range("N132").Value = -23058.9995
$ws.Range("H134").Value = 2544.5
$ws.Range("I134").Value = 1299.909
$ws.Range("J134").Value = 3789.0908
$ws.Range("K134").Value = 3899.727
$ws.Range("L134").Value = 11367.2724
$ws.Range("M134").Value = -1364.727
$ws.Range("N134").Value = -16437.2724
$ws.Range("H136").Value = 1503.086
$ws.Range("I136").Value = 1268.6666
$ws.Range("J136").Value = 3691
$ws.Range("K136").Value = 3805.9998
$ws.Range("L136").Value = 11073
$ws.Range("M136").Value = -1255.9998
$ws.Range("N136").Value = -16173

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 5039.3335
$ws.Range("I109").Value = 4451.3335
$ws.Range("J109").Value = 5333.3335
$ws.Range("K109").Value = 13354.0005
$ws.Range("L109").Value = 16000.0005
$ws.Range("M109").Value = -12314.0005
$ws.Range("N109").Value = -18080.0005
$ws.Range("H131").Value = 7143657
$ws.Range("J131").Value = 871.39343
$ws.Range("L131").Value = 2614.18029
$ws.Range("N131").Value = -12694.18029
$ws.Range("H137").Value = 2463.647
$ws.Range("I137").Value = 665
$ws.Range("J137").Value = 5033.143
$ws.Range("K137").Value = 1995
$ws.Range("L137").Value = 15099.429
$ws.Range("M137").Value = 3105
$ws.Range("N137").Value = -25299.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2214.3125
$ws.Range("I132").Value = 1270.5135
$ws.Range("J132").Value = 5388.909
$ws.Range("K132").Value = 3811.5405
$ws.Range("L132").Value = 16166.727
$ws.Range("M132").Value = -1281.5405
$ws.Range("N132").Value = -21226.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2440
$ws.Range("I100").Value = 2066.6667
$ws.Range("K100").Value = 2066.6667
$ws.Range("M100").Value = -1525.6667
$ws.Range("H132").Value = 10911.027
$ws.Range("I132").Value = 11050.154
$ws.Range("J132").Value = 10582.182
$ws.Range("K132").Value = 33150.462
$ws.Range("L132").Value = 31746.546
$ws.Range("M132").Value = -30620.462
$ws.Range("N132").Value = -36806.546
$ws.Range("H136").Value = 2087.35
$ws.Range("I136").Value = 1151.9395
$ws.Range("K136").Value = 3455.8185
$ws.Range("M136").Value = -905.8184999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H131").Value = 52978.57
$ws.Range("J131").Value = 52978.57
$ws.Range("L131").Value = 52978.57
$ws.Range("N131").Value = -63058.57
$ws.Range("H132").Value = 4903472
$ws.Range("I132").Value = 1099.0212
$ws.Range("J132").Value = 15875450
$ws.Range("K132").Value = 3297.063599999999
$ws.Range("L132").Value = 47626350
$ws.Range("M132").Value = -767.0635999999995
$ws.Range("N132").Value = -47631410
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()
